$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(289).Insert()

$ws.Cells.Item(289, 1).Value = 4
$ws.Cells.Item(289, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(289, 3).Value = "Los Lagos"
$ws.Cells.Item(289, 4).Value = 44868
$ws.Cells.Item(289, 5).Value = 10
$ws.Cells.Item(289, 6).Value = "Fruta"
$ws.Cells.Item(289, 7).Value = 100104
$ws.Cells.Item(289, 8).Value = "Frutos de pepita"
$ws.Cells.Item(289, 9).Value = 100104005
$ws.Cells.Item(289, 10).Value = "Pera"
$ws.Cells.Item(289, 11).Value = "Packham's Triumph"
$ws.Cells.Item(289, 12).Value = "Especial"
$ws.Cells.Item(289, 13).Value = 150
$ws.Cells.Item(289, 14).Value = 22000
$ws.Cells.Item(289, 15).Value = 22000
$ws.Cells.Item(289, 16).Value = 22000
$ws.Cells.Item(289, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(289, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(289, 19).Value = 1467
$ws.Cells.Item(289, 20).Value = 15
